$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.162.80'
$ws.Range("E2").Value = '  -1.25%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.657.27'
$ws.Range("E3").Value = '  -1.25%  '
$ws.Range("E4").Value = '  +0.45%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.66'
$ws.Range("E5").Value = '  -1.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5186'
$ws.Range("E6").Value = '  -2.89%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2624'
$ws.Range("E8").Value = '  -2.77%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06232'
$ws.Range("E9").Value = '  -2.90%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '20.61'
$ws.Range("E10").Value = '  -6.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07707'
$ws.Range("E11").Value = '  -1.13%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.667.83'
$ws.Range("E12").Value = '  -0.57%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.389'
$ws.Range("E13").Value = '  -2.84%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.884.65'
$ws.Range("E14").Value = '  -1.20%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5403'
$ws.Range("E15").Value = '  -3.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅8068'
$ws.Range("E16").Value = '  -3.25%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.14'
$ws.Range("E17").Value = '  -2.53%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '26.189.50'
$ws.Range("E18").Value = '  -1.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.597'
$ws.Range("E20").Value = '  -4.04%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '190.34'
$ws.Range("E21").Value = '  -1.81%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.991'
$ws.Range("E22").Value = '  -3.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.035'
$ws.Range("E23").Value = '  -4.88%  '
$ws.Range("E24").Value = '  +0.55%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '140.01'
$ws.Range("E25").Value = '  -1.67%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1218'
$ws.Range("E26").Value = '  -5.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.116'
$ws.Range("E27").Value = '  -4.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.93'
$ws.Range("E28").Value = '  -2.60%  '
$ws.Range("E29").Value = '  -2.59%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.05962'
$ws.Range("E30").Value = '  -5.04%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.266'
$ws.Range("E31").Value = '  -0.76%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.551'
$ws.Range("E32").Value = '  -1.80%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.230'
$ws.Range("E33").Value = '  -6.83%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.611'
$ws.Range("E34").Value = '  -5.36%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9596'
$ws.Range("E35").Value = '  -5.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.418'
$ws.Range("E36").Value = '  -0.13%  '
$ws.Range("E37").Value = '  -0.43%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5646'
$ws.Range("E38").Value = '  -7.70%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.000'
$ws.Range("E39").Value = '  -2.83%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01590'
$ws.Range("E40").Value = '  -2.66%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8556'
$ws.Range("E41").Value = '  -0.86%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.004'
$ws.Range("E42").Value = '  +0.46%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.016.07'
$ws.Range("E43").Value = '  -7.27%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.90'
$ws.Range("E44").Value = '  -0.68%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.799.92'
$ws.Range("E46").Value = '  +2.57%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '56.46'
$ws.Range("E47").Value = '  -2.57%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.013'
$ws.Range("E48").Value = '  +0.90%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.944'
$ws.Range("E49").Value = '  -2.72%  '
$ws.Range("E50").Value = '  -0.62%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4209'
$ws.Range("E51").Value = '  -0.60%  '
